# Auto-generated edit script: apply numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 81
$ws.Range("H81").Value = 34204
$ws.Range("J81").Value = 34204
$ws.Range("L81").Value = 34204
$ws.Range("N81").Value = -36200

# Row 84
$ws.Range("H84").Value = 34204
$ws.Range("J84").Value = 34204
$ws.Range("L84").Value = 102612
$ws.Range("N84").Value = -112596

# Row 137
$ws.Range("H137").Value = 14057149
$ws.Range("I137").Value = 2488387
$ws.Range("J137").Value = 62501340
$ws.Range("K137").Value = 7465161
$ws.Range("L137").Value = 187504020
$ws.Range("M137").Value = -7462611
$ws.Range("N137").Value = -187509120

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 13167.667
$ws.Range("J8").Value = 13167.667
$ws.Range("L8").Value = 13167.667
$ws.Range("N8").Value = -13455.667

# Row 32
$ws.Range("H32").Value = 3918.045
$ws.Range("I32").Value = 3914.4236
$ws.Range("J32").Value = 3995
$ws.Range("K32").Value = 3914.4236
$ws.Range("L32").Value = 3995
$ws.Range("M32").Value = -3627.4236
$ws.Range("N32").Value = -4569

# Row 61
$ws.Range("H61").Value = 5953522
$ws.Range("I61").Value = 7247473.5
$ws.Range("J61").Value = 1345.6
$ws.Range("K61").Value = 7247473.5
$ws.Range("L61").Value = 1345.6
$ws.Range("M61").Value = -7247261.5
$ws.Range("N61").Value = -1769.6

# Row 74
$ws.Range("H74").Value = 7696688.5
$ws.Range("I74").Value = 12000671
$ws.Range("J74").Value = 11006.857
$ws.Range("K74").Value = 12000671
$ws.Range("L74").Value = 11006.857
$ws.Range("M74").Value = -11999797
$ws.Range("N74").Value = -12754.857

# Row 77
$ws.Range("H77").Value = 7696688.5
$ws.Range("I77").Value = 12000671
$ws.Range("J77").Value = 11006.857
$ws.Range("K77").Value = 60003355
$ws.Range("L77").Value = 55034.285
$ws.Range("M77").Value = -59998987
$ws.Range("N77").Value = -63770.285

# Row 136
$ws.Range("H136").Value = 5953522
$ws.Range("I136").Value = 7247473.5
$ws.Range("J136").Value = 1345.6
$ws.Range("K136").Value = 21742420.5
$ws.Range("L136").Value = 4036.8
$ws.Range("M136").Value = -21739870.5
$ws.Range("N136").Value = -9136.799999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 39056.723
$ws.Range("I134").Value = 56064.45
$ws.Range("J134").Value = 1261.7778
$ws.Range("K134").Value = 168193.35
$ws.Range("L134").Value = 3785.3334
$ws.Range("M134").Value = -165658.35
$ws.Range("N134").Value = -8855.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Row 31
$ws.Range("H31").Value = 1682.2
$ws.Range("I31").Value = 963.2692
$ws.Range("J31").Value = 2666
$ws.Range("K31").Value = 963.2692
$ws.Range("L31").Value = 2666
$ws.Range("M31").Value = -668.2692
$ws.Range("N31").Value = -3256

# Row 34
$ws.Range("H34").Value = 1682.2
$ws.Range("I34").Value = 963.2692
$ws.Range("J34").Value = 2666
$ws.Range("K34").Value = 963.2692
$ws.Range("L34").Value = 2666
$ws.Range("M34").Value = -761.2692
$ws.Range("N34").Value = -3070

# Row 58
$ws.Range("H58").Value = 948.7925
$ws.Range("I58").Value = 766.913
$ws.Range("J58").Value = 2144
$ws.Range("K58").Value = 766.913
$ws.Range("L58").Value = 2144
$ws.Range("M58").Value = -563.913
$ws.Range("N58").Value = -2550

# Row 107
$ws.Range("H107").Value = 313.52173
$ws.Range("I107").Value = 251.3125
$ws.Range("J107").Value = 455.7143
$ws.Range("K107").Value = 251.3125
$ws.Range("L107").Value = 455.7143
$ws.Range("M107").Value = 1668.6875
$ws.Range("N107").Value = -4295.7143

# Row 132
$ws.Range("H132").Value = 2064.7424
$ws.Range("I132").Value = 1968.0518
$ws.Range("J132").Value = 2765.75
$ws.Range("K132").Value = 5904.1554
$ws.Range("L132").Value = 8297.25
$ws.Range("M132").Value = -3374.1554
$ws.Range("N132").Value = -13357.25

# Row 134
$ws.Range("H134").Value = 2189
$ws.Range("I134").Value = 2408.353
$ws.Range("J134").Value = 1256.75
$ws.Range("K134").Value = 7225.059
$ws.Range("L134").Value = 3770.25
$ws.Range("M134").Value = -4690.059
$ws.Range("N134").Value = -8840.25

# Row 136
$ws.Range("H136").Value = 948.7925
$ws.Range("I136").Value = 766.913
$ws.Range("J136").Value = 2144
$ws.Range("K136").Value = 2300.739
$ws.Range("L136").Value = 6432
$ws.Range("M136").Value = 249.261
$ws.Range("N136").Value = -11532

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 1973.75
$ws.Range("I3").Value = 947.5
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 2842.5
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -2730.5
$ws.Range("N3").Value = -9224

# Row 98
$ws.Range("H98").Value = 564.4286
$ws.Range("I98").Value = 115.666664
$ws.Range("J98").Value = 901
$ws.Range("K98").Value = 346.999992
$ws.Range("L98").Value = 2703
$ws.Range("M98").Value = 1151.000008
$ws.Range("N98").Value = -5699

# Row 125
$ws.Range("H125").Value = 2293.6365
$ws.Range("I125").Value = 1686
$ws.Range("J125").Value = 2800
$ws.Range("K125").Value = 5058
$ws.Range("L125").Value = 8400
$ws.Range("M125").Value = -138
$ws.Range("N125").Value = -18240

$ws = $wb.Worksheets.Item("GSM")
# Row 125
$ws.Range("H125").Value = 36000
$ws.Range("J125").Value = 36000
$ws.Range("L125").Value = 36000
$ws.Range("N125").Value = -40920

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 547
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 444
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 444
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1034

# Row 27
$ws.Range("H27").Value = 547
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 444
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 444
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -658

# Row 132
$ws.Range("H132").Value = 28087.44
$ws.Range("I132").Value = 28639.65
$ws.Range("K132").Value = 85918.95000000001
$ws.Range("M132").Value = -83388.95000000001

# Row 136
$ws.Range("H136").Value = 4712.75
$ws.Range("I136").Value = 4559.885
$ws.Range("J136").Value = 6700
$ws.Range("K136").Value = 13679.655
$ws.Range("L136").Value = 20100
$ws.Range("M136").Value = -11129.655
$ws.Range("N136").Value = -25200

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 7501577
$ws.Range("I3").Value = 15000154
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 15000154
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -15000040
$ws.Range("N3").Value = -3228

# Row 11
$ws.Range("H11").Value = 50000000
$ws.Range("I11").Value = 50000000
$ws.Range("K11").Value = 50000000
$ws.Range("M11").Value = -49999858

# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

# Row 132
$ws.Range("H132").Value = 9479.842000000001
$ws.Range("I132").Value = 10557.4375
$ws.Range("J132").Value = 3732.6667
$ws.Range("K132").Value = 31672.3125
$ws.Range("L132").Value = 11198.0001
$ws.Range("M132").Value = -29142.3125
$ws.Range("N132").Value = -16258.0001

# Row 136
$ws.Range("H136").Value = 12626.294
$ws.Range("I136").Value = 14109.8
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 42329.39999999999
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -39779.39999999999
$ws.Range("N136").Value = -9600

